$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the two new example sheets ("skiprow" and "usecols") at the end of
# the workbook, after the existing "colheader" sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$skiprow = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$skiprow.Name = "skiprow"

$usecols = $wb.Worksheets.Add([System.Type]::Missing, $skiprow)
$usecols.Name = "usecols"

# ---------------------------------------------------------------------
# "skiprow" sheet data - a header row, a couple of junk "no" rows, then
# the real data (used to illustrate pandas read_excel(skiprows=...)).
# ---------------------------------------------------------------------
$skiprow.Range("A1").Value = "No"
$skiprow.Range("B1").Value = "city"

$skiprow.Range("A2").Value = 1
$skiprow.Range("B2").Value = "no"

$skiprow.Range("A3").Value = 2
$skiprow.Range("B3").Value = "no"

$skiprow.Range("A4").Value = 3
$skiprow.Range("B4").Value = "no"

$skiprow.Range("A5").Value = 4
$skiprow.Range("B5").Value = "Seoul"

$skiprow.Range("A6").Value = 5
$skiprow.Range("B6").Value = "Pusan"

$skiprow.Range("A1:B6").Select() | Out-Null

$skiprow.PageSetup.PaperSize = 9
$skiprow.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# "usecols" sheet data - two side-by-side tables (A:B and D:F), used to
# illustrate pandas read_excel(usecols=...).
# ---------------------------------------------------------------------
$usecols.Range("A1").Value = "No"
$usecols.Range("B1").Value = "city"
$usecols.Range("D1").Value = "No"
$usecols.Range("E1").Value = "city"
$usecols.Range("F1").Value = "population"

$usecols.Range("A2").Value = 0
$usecols.Range("B2").Value = "Seoul"
$usecols.Range("D2").Value = 1
$usecols.Range("E2").Value = "Seoul"
$usecols.Range("F2").Value = 10000000

$usecols.Range("A3").Value = 2
$usecols.Range("B3").Value = "Pusan"
$usecols.Range("D3").Value = 2
$usecols.Range("E3").Value = "Pusan"
$usecols.Range("F3").Value = 5000000

$usecols.Columns.Item(6).ColumnWidth = 8.36

# Active cell / selection + which sheet/tab is active when the workbook
# is opened.
$usecols.Range("A2").Select() | Out-Null
